$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume Number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 32   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/19/2025  Through  5/25/2025"

# --- Row 15 ---
$ws.Range("N15").Value = -66.666666666666

# --- Row 16 ---
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 12
$ws.Range("H16").Value = 9.090909090909
$ws.Range("I16").Value = 70
$ws.Range("J16").Value = 73
$ws.Range("K16").Value = -4.109589041095
$ws.Range("L16").Value = -7.894736842105
$ws.Range("M16").Value = 37.254901960784
$ws.Range("N16").Value = -86.86679174484

# --- Row 17 ---
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -27.272727272727
$ws.Range("I17").Value = 75
$ws.Range("J17").Value = 71
$ws.Range("K17").Value = 5.633802816901
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 102.702702702703
$ws.Range("N17").Value = -33.62831858407

# --- Row 18 ---
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 66.666666666666
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = 15.78947368421
$ws.Range("I18").Value = 119
$ws.Range("J18").Value = 98
$ws.Range("K18").Value = 21.428571428571
$ws.Range("L18").Value = 1.709401709401
$ws.Range("M18").Value = 11.214953271028
$ws.Range("N18").Value = -88.741721854304

# --- Row 19 ---
$ws.Range("C19").Value = 45
$ws.Range("D19").Value = 32
$ws.Range("E19").Value = 40.625
$ws.Range("F19").Value = 139
$ws.Range("G19").Value = 120
$ws.Range("H19").Value = 15.833333333333
$ws.Range("I19").Value = 619
$ws.Range("J19").Value = 608
$ws.Range("K19").Value = 1.809210526315
$ws.Range("L19").Value = -2.211690363349
$ws.Range("M19").Value = 29.227557411273
$ws.Range("N19").Value = -58.203916272788

# --- Row 20 ---
$ws.Range("C20").NumberFormat = '#,##0'
$ws.Range("C20").Value = 3
$ws.Range("D20").NumberFormat = '#,##0'
$ws.Range("D20").Value = 4
$ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 21
$ws.Range("J20").Value = 25
$ws.Range("K20").Value = -16
$ws.Range("L20").Value = -67.1875
$ws.Range("M20").Value = -12.5
$ws.Range("N20").Value = -98.498927805575

# --- Row 21 ---
$ws.Range("C21").Value = 60
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = 36.363636363636
$ws.Range("F21").Value = 195
$ws.Range("G21").Value = 178
$ws.Range("H21").Value = 9.550561797752
$ws.Range("I21").Value = 907
$ws.Range("J21").Value = 881
$ws.Range("K21").Value = 2.951191827468
$ws.Range("L21").Value = -6.687242798353
$ws.Range("M21").Value = 28.835227272727
$ws.Range("N21").Value = -80.261153427638

# --- Row 22 ---
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 2
$ws.Range("D22").NumberFormat = 'General'
$ws.Range("D22").Value = "0"
$ws.Range("E22").NumberFormat = 'General'
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 14
$ws.Range("K22").Value = -12.5
$ws.Range("L22").Value = -26.315789473684
$ws.Range("M22").Value = 27.272727272727

# --- Row 23 ---
$ws.Range("D23").NumberFormat = '#,##0'
$ws.Range("D23").Value = 1
$ws.Range("E23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E23").Value = -100
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = -50
$ws.Range("J23").Value = 10
$ws.Range("K23").Value = 30

# --- Row 24 ---
$ws.Range("C24").Value = 42
$ws.Range("E24").Value = -27.586206896551
$ws.Range("F24").Value = 181
$ws.Range("G24").Value = 219
$ws.Range("H24").Value = -17.351598173516
$ws.Range("I24").Value = 1213
$ws.Range("J24").Value = 1136
$ws.Range("K24").Value = 6.778169014084
$ws.Range("L24").Value = -0.164609053497
$ws.Range("M24").Value = 93.460925039872

# --- Row 25 ---
$ws.Range("C25").Value = 37
$ws.Range("D25").Value = 50
$ws.Range("E25").Value = -26
$ws.Range("F25").Value = 148
$ws.Range("G25").Value = 197
$ws.Range("H25").Value = -24.8730964467
$ws.Range("I25").Value = 1033
$ws.Range("J25").Value = 1002
$ws.Range("K25").Value = 3.093812375249
$ws.Range("L25").Value = -2.913533834586

# --- Row 26 ---
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 15
$ws.Range("E26").Value = -40
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = -17.948717948717
$ws.Range("I26").Value = 128
$ws.Range("J26").Value = 144
$ws.Range("K26").Value = -11.111111111111
$ws.Range("L26").Value = -6.569343065693
$ws.Range("M26").Value = -11.724137931034

# --- Row 27 ---
$ws.Range("L27").Value = -50

# --- Row 28 ---
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = -75
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = -38.461538461538
$ws.Range("I28").Value = 34
$ws.Range("J28").Value = 56
$ws.Range("K28").Value = -39.285714285714
$ws.Range("L28").Value = -19.047619047619

# --- Row 31 ---
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = -66.666666666666
$ws.Range("J31").Value = 14
$ws.Range("K31").Value = -42.857142857142
